$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date values in column C
$ws.Range("C3").Value = 44926
$ws.Range("C4").Value = 45291
$ws.Range("C6").Value = 44926
$ws.Range("C7").Value = 45291
$ws.Range("C9").Value = 44926
$ws.Range("C10").Value = 45291

# Update the selected range on the sheet
$ws.Range("C8:C10").Select()
